# Updated cryptos list on Tue Sep 19 05:37:36 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto ranking sheet, and fixes the MXToken/Aave rows (44-45) which had
# been swapped - each coin's name/link/price/volume now line up correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value still "looks like a number" (e.g. "19.88") need to
# be forced to Text format first, otherwise Excel auto-converts the string
# into a numeric value when it is assigned via .Value.
$numberLikeCells = @(
    "D5","D10","D15","D16","D19","D21","D25","D29","D32","D39","D42",
    "D44","D45","D46","D47","D51"
)
foreach ($cellRef in $numberLikeCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "26.873.86"
$ws.Range("E2").Value = "  +0.21%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "1.639.83"
$ws.Range("E3").Value = "  -0.24%  "

# --- Row 5: BNB ---
$ws.Range("D5").Value = "216.82"
$ws.Range("E5").Value = "  -0.77%  "

# --- Row 6: XRP ---
$ws.Range("E6").Value = "  +1.84%  "

# --- Row 8: Cardano ---
$ws.Range("E8").Value = "  +1.75%  "

# --- Row 9: Dogecoin ---
$ws.Range("E9").Value = "  +0.43%  "

# --- Row 10: Solana ---
$ws.Range("D10").Value = "19.88"
$ws.Range("E10").Value = "  +3.35%  "

# --- Row 11: TRON ---
$ws.Range("E11").Value = "  -0.07%  "

# --- Row 12: WrappedliquidstakedEther2.0 ---
$ws.Range("D12").Value = "1.869.17"
$ws.Range("E12").Value = "  -0.18%  "

# --- Row 13: WrappedEther ---
$ws.Range("D13").Value = "1.638.46"
$ws.Range("E13").Value = "  -0.01%  "

# --- Row 14: Polkadot ---
$ws.Range("E14").Value = "  -0.71%  "

# --- Row 15: Polygon ---
$ws.Range("D15").Value = "0.531"
$ws.Range("E15").Value = "  +0.92%  "

# --- Row 16: Litecoin ---
$ws.Range("D16").Value = "67.29"
$ws.Range("E16").Value = "  +3.11%  "

# --- Row 17: WrappedBTC ---
$ws.Range("D17").Value = "26.868.04"
$ws.Range("E17").Value = "  +0.18%  "

# --- Row 18: ShibaInu ---
$ws.Range("E18").Value = "  -0.64%  "

# --- Row 19: BitcoinCash ---
$ws.Range("D19").Value = "219.90"
$ws.Range("E19").Value = "  +2.07%  "

# --- Row 20: Dai ---
$ws.Range("E20").Value = "  -0.53%  "

# --- Row 21: Chainlink ---
$ws.Range("D21").Value = "6.86"
$ws.Range("E21").Value = "  +3.69%  "

# --- Row 22: Uniswap ---
$ws.Range("E22").Value = "  +0.41%  "

# --- Row 23: Toncoin ---
$ws.Range("E23").Value = "  +3.68%  "

# --- Row 24: Avalanche ---
$ws.Range("E24").Value = "  -0.47%  "

# --- Row 25: Monero ---
$ws.Range("D25").Value = "147.15"
$ws.Range("E25").Value = "  -0.40%  "

# --- Row 26: BinanceUSD ---
$ws.Range("E26").Value = "  -0.54%  "

# --- Row 27: Cosmos (etc.) ---
$ws.Range("E27").Value = "  +2.98%  "

# --- Row 28 ---
$ws.Range("E28").Value = "  +0.28%  "

# --- Row 29 ---
$ws.Range("D29").Value = "15.80"
$ws.Range("E29").Value = "  +0.48%  "

# --- Row 30 ---
$ws.Range("E30").Value = "  -1.21%  "

# --- Row 31 ---
$ws.Range("E31").Value = "  -0.87%  "

# --- Row 32 ---
$ws.Range("D32").Value = "3.34"
$ws.Range("E32").Value = "  -1.32%  "

# --- Row 33 ---
$ws.Range("E33").Value = "  +0.50%  "

# --- Row 34 ---
$ws.Range("E34").Value = "  +1.32%  "

# --- Row 35 ---
$ws.Range("D35").Value = "1.265.50"
$ws.Range("E35").Value = "  -0.26%  "

# --- Row 36 ---
$ws.Range("E36").Value = "  -0.14%  "

# --- Row 37 ---
$ws.Range("E37").Value = "  +2.06%  "

# --- Row 38 ---
$ws.Range("E38").Value = "  +0.34%  "

# --- Row 39 ---
$ws.Range("D39").Value = "0.834"
$ws.Range("E39").Value = "  +2.08%  "

# --- Row 40 ---
$ws.Range("E40").Value = "  -0.45%  "

# --- Row 41 ---
$ws.Range("E41").Value = "  +0.98%  "

# --- Row 42: FraxShare ---
$ws.Range("D42").Value = "5.39"

# --- Row 43 ---
$ws.Range("D43").Value = "1.779.39"
$ws.Range("E43").Value = "  -0.13%  "

# --- Row 44: was MXToken, now Aave ---
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "61.87"
$ws.Range("E44").Value = "  +0.70%  "

# --- Row 45: was Aave, now MXToken ---
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").Value = "2.11"
$ws.Range("E45").Value = "  -1.48%  "

# --- Row 46 ---
$ws.Range("D46").Value = "91.82"
$ws.Range("E46").Value = "  -1.10%  "

# --- Row 47 ---
$ws.Range("D47").Value = "1.59"
$ws.Range("E47").Value = "  -0.97%  "

# --- Row 48 ---
$ws.Range("E48").Value = "  +1.12%  "

# --- Row 49 ---
$ws.Range("E49").Value = "  -0.42%  "

# --- Row 50 ---
$ws.Range("E50").Value = "  +1.52%  "

# --- Row 51 ---
$ws.Range("D51").Value = "0.0962"
$ws.Range("E51").Value = "  -0.31%  "
